$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Hydrogen): update B3 value, clear D3 (becomes blank / inlineStr)
$ws.Range("B3").Value = 432144.236637009
$ws.Range("D3").ClearContents()

# Row 4 (Methanol): update C4 value
$ws.Range("C4").Value = 30.49981016068243

# Row 5 (Ammonia): update C5 value
$ws.Range("C5").Value = 0

# Row 7: rename "Other" -> "Biogas", update D7 value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 35.7588053685814

# New row 8: "Other" entry, copying the formatting of row 7's label cell
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "Other"
$ws.Range("D8").Value = 179.7972293456137
